$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 6 new rows at the top of the data (rows 2-7) for the new
#    "05/2025" (May/2025) month, and fill them with the new daily data.
# ---------------------------------------------------------------------------
$ws.Range("A2:E7").EntireRow.Insert()
$ws.Range("A2:E7").ClearFormats()

$mayData = @(
    @(1, 17056.87),
    @(2, 22786.63),
    @(3, 29255),
    @(4, 4370),
    @(5, 35402.64),
    @(6, 17735.42)
)

$r = 2
foreach ($row in $mayData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = 5
    $ws.Cells.Item($r, 4).Value = 2025
    $ws.Cells.Item($r, 5).Value = "05/2025"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Add a new row for 04/2025 (April) day 30, which was previously missing.
#    After the insert above, April's data now occupies rows 8-36, so the new
#    day 30 entry goes right after, at row 37.
# ---------------------------------------------------------------------------
$ws.Range("A37:E37").EntireRow.Insert()
$ws.Range("A37:E37").ClearFormats()
$ws.Cells.Item(37, 1).Value = 30
$ws.Cells.Item(37, 2).Value = 24773.88
$ws.Cells.Item(37, 3).Value = 4
$ws.Cells.Item(37, 4).Value = 2025
$ws.Cells.Item(37, 5).Value = "04/2025"

# ---------------------------------------------------------------------------
# 3) Remove the oldest month, 01/2025 (January), entirely. After the two
#    inserts above its rows shifted from 90-119 to 97-126.
# ---------------------------------------------------------------------------
$ws.Range("A97:E126").EntireRow.Delete()
